$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target column widths (as they should appear in the saved OOXML <col width="..."/>)
# Excel's ColumnWidth property (in characters) differs from the stored width by
# a constant offset of 5/6 (0.8333...) due to internal padding, so subtract it
# here to land exactly on the desired stored widths: 23,19,20,20,20,20,20,19,20
$offset = 5 / 6

$ws.Columns.Item(1).ColumnWidth = 23 - $offset
$ws.Columns.Item(2).ColumnWidth = 19 - $offset
$ws.Columns.Item(3).ColumnWidth = 20 - $offset
$ws.Columns.Item(4).ColumnWidth = 20 - $offset
$ws.Columns.Item(5).ColumnWidth = 20 - $offset
$ws.Columns.Item(6).ColumnWidth = 20 - $offset
$ws.Columns.Item(7).ColumnWidth = 20 - $offset
$ws.Columns.Item(8).ColumnWidth = 19 - $offset
$ws.Columns.Item(9).ColumnWidth = 20 - $offset
